$wb = $excel.ActiveWorkbook

# Data-driven cell updates derived from the authoritative diff.
# Each entry: worksheet name, cell reference, and the new numeric value.
$edits = @(
    @{ Sheet = "ALC"; Cell = "H28"; Value = 1712.4445 },
    @{ Sheet = "ALC"; Cell = "I28"; Value = 1300 },
    @{ Sheet = "ALC"; Cell = "J28"; Value = 1918.6666 },
    @{ Sheet = "ALC"; Cell = "K28"; Value = 1300 },
    @{ Sheet = "ALC"; Cell = "L28"; Value = 1918.6666 },
    @{ Sheet = "ALC"; Cell = "M28"; Value = -815 },
    @{ Sheet = "ALC"; Cell = "N28"; Value = -2888.6666 },
    @{ Sheet = "ALC"; Cell = "H33"; Value = 136.83871 },
    @{ Sheet = "ALC"; Cell = "I33"; Value = 91.07692 },
    @{ Sheet = "ALC"; Cell = "K33"; Value = 91.07692 },
    @{ Sheet = "ALC"; Cell = "M33"; Value = 137.92308 },
    @{ Sheet = "ALC"; Cell = "H74"; Value = 5060.375 },
    @{ Sheet = "ALC"; Cell = "I74"; Value = 4580.5 },
    @{ Sheet = "ALC"; Cell = "J74"; Value = 6500 },
    @{ Sheet = "ALC"; Cell = "K74"; Value = 4580.5 },
    @{ Sheet = "ALC"; Cell = "L74"; Value = 6500 },
    @{ Sheet = "ALC"; Cell = "M74"; Value = -3644.5 },
    @{ Sheet = "ALC"; Cell = "N74"; Value = -8372 },
    @{ Sheet = "ALC"; Cell = "H77"; Value = 5060.375 },
    @{ Sheet = "ALC"; Cell = "I77"; Value = 4580.5 },
    @{ Sheet = "ALC"; Cell = "J77"; Value = 6500 },
    @{ Sheet = "ALC"; Cell = "K77"; Value = 22902.5 },
    @{ Sheet = "ALC"; Cell = "L77"; Value = 32500 },
    @{ Sheet = "ALC"; Cell = "M77"; Value = -18222.5 },
    @{ Sheet = "ALC"; Cell = "N77"; Value = -41860 },
    @{ Sheet = "ALC"; Cell = "H112"; Value = 1537.3125 },
    @{ Sheet = "ALC"; Cell = "J112"; Value = 1613.1333 },
    @{ Sheet = "ALC"; Cell = "L112"; Value = 4839.3999 },
    @{ Sheet = "ALC"; Cell = "N112"; Value = -7055.3999 },
    @{ Sheet = "ALC"; Cell = "H125"; Value = 8229.308000000001 },
    @{ Sheet = "ALC"; Cell = "I125"; Value = 450 },
    @{ Sheet = "ALC"; Cell = "J125"; Value = 9643.727999999999 },
    @{ Sheet = "ALC"; Cell = "K125"; Value = 4050 },
    @{ Sheet = "ALC"; Cell = "L125"; Value = 86793.552 },
    @{ Sheet = "ALC"; Cell = "M125"; Value = -1590 },
    @{ Sheet = "ALC"; Cell = "N125"; Value = -91713.552 },
    @{ Sheet = "ARM"; Cell = "H2"; Value = 2516.7856 },
    @{ Sheet = "ARM"; Cell = "I2"; Value = 2935.7778 },
    @{ Sheet = "ARM"; Cell = "J2"; Value = 1762.6 },
    @{ Sheet = "ARM"; Cell = "K2"; Value = 2935.7778 },
    @{ Sheet = "ARM"; Cell = "L2"; Value = 1762.6 },
    @{ Sheet = "ARM"; Cell = "M2"; Value = -2822.7778 },
    @{ Sheet = "ARM"; Cell = "N2"; Value = -1988.6 },
    @{ Sheet = "ARM"; Cell = "H32"; Value = 18350.73 },
    @{ Sheet = "ARM"; Cell = "I32"; Value = 20531.277 },
    @{ Sheet = "ARM"; Cell = "K32"; Value = 20531.277 },
    @{ Sheet = "ARM"; Cell = "M32"; Value = -20244.277 },
    @{ Sheet = "ARM"; Cell = "H63"; Value = 4150.7144 },
    @{ Sheet = "ARM"; Cell = "I63"; Value = 3513.75 },
    @{ Sheet = "ARM"; Cell = "J63"; Value = 5000 },
    @{ Sheet = "ARM"; Cell = "K63"; Value = 3513.75 },
    @{ Sheet = "ARM"; Cell = "L63"; Value = 5000 },
    @{ Sheet = "ARM"; Cell = "M63"; Value = -2827.75 },
    @{ Sheet = "ARM"; Cell = "N63"; Value = -6372 },
    @{ Sheet = "ARM"; Cell = "H66"; Value = 4150.7144 },
    @{ Sheet = "ARM"; Cell = "I66"; Value = 3513.75 },
    @{ Sheet = "ARM"; Cell = "J66"; Value = 5000 },
    @{ Sheet = "ARM"; Cell = "K66"; Value = 17568.75 },
    @{ Sheet = "ARM"; Cell = "L66"; Value = 25000 },
    @{ Sheet = "ARM"; Cell = "M66"; Value = -14136.75 },
    @{ Sheet = "ARM"; Cell = "N66"; Value = -31864 },
    @{ Sheet = "ARM"; Cell = "H74"; Value = 1626.5814 },
    @{ Sheet = "ARM"; Cell = "I74"; Value = 1636.0526 },
    @{ Sheet = "ARM"; Cell = "K74"; Value = 1636.0526 },
    @{ Sheet = "ARM"; Cell = "M74"; Value = -762.0526 },
    @{ Sheet = "ARM"; Cell = "H77"; Value = 1626.5814 },
    @{ Sheet = "ARM"; Cell = "I77"; Value = 1636.0526 },
    @{ Sheet = "ARM"; Cell = "K77"; Value = 8180.263 },
    @{ Sheet = "ARM"; Cell = "M77"; Value = -3812.263 },
    @{ Sheet = "ARM"; Cell = "H116"; Value = 2516.7856 },
    @{ Sheet = "ARM"; Cell = "I116"; Value = 2935.7778 },
    @{ Sheet = "ARM"; Cell = "J116"; Value = 1762.6 },
    @{ Sheet = "ARM"; Cell = "K116"; Value = 2935.7778 },
    @{ Sheet = "ARM"; Cell = "L116"; Value = 1762.6 },
    @{ Sheet = "ARM"; Cell = "M116"; Value = -641.7777999999998 },
    @{ Sheet = "ARM"; Cell = "N116"; Value = -6350.6 },
    @{ Sheet = "ARM"; Cell = "H122"; Value = 2800 },
    @{ Sheet = "ARM"; Cell = "I122"; Value = 2700 },
    @{ Sheet = "ARM"; Cell = "J122"; Value = 2950 },
    @{ Sheet = "ARM"; Cell = "K122"; Value = 8100 },
    @{ Sheet = "ARM"; Cell = "L122"; Value = 8850 },
    @{ Sheet = "ARM"; Cell = "M122"; Value = -5650 },
    @{ Sheet = "ARM"; Cell = "N122"; Value = -13750 },
    @{ Sheet = "BSM"; Cell = "H3"; Value = 2516.7856 },
    @{ Sheet = "BSM"; Cell = "I3"; Value = 2935.7778 },
    @{ Sheet = "BSM"; Cell = "J3"; Value = 1762.6 },
    @{ Sheet = "BSM"; Cell = "K3"; Value = 2935.7778 },
    @{ Sheet = "BSM"; Cell = "L3"; Value = 1762.6 },
    @{ Sheet = "BSM"; Cell = "M3"; Value = -2821.7778 },
    @{ Sheet = "BSM"; Cell = "N3"; Value = -1990.6 },
    @{ Sheet = "BSM"; Cell = "H82"; Value = 36529.5 },
    @{ Sheet = "BSM"; Cell = "I82"; Value = 10257 },
    @{ Sheet = "BSM"; Cell = "J82"; Value = 40282.715 },
    @{ Sheet = "BSM"; Cell = "K82"; Value = 10257 },
    @{ Sheet = "BSM"; Cell = "L82"; Value = 40282.715 },
    @{ Sheet = "BSM"; Cell = "M82"; Value = -9874 },
    @{ Sheet = "BSM"; Cell = "N82"; Value = -41048.715 },
    @{ Sheet = "BSM"; Cell = "H85"; Value = 36529.5 },
    @{ Sheet = "BSM"; Cell = "I85"; Value = 10257 },
    @{ Sheet = "BSM"; Cell = "J85"; Value = 40282.715 },
    @{ Sheet = "BSM"; Cell = "K85"; Value = 10257 },
    @{ Sheet = "BSM"; Cell = "L85"; Value = 40282.715 },
    @{ Sheet = "BSM"; Cell = "M85"; Value = -8931 },
    @{ Sheet = "BSM"; Cell = "N85"; Value = -42934.715 },
    @{ Sheet = "CRP"; Cell = "H99"; Value = 1217.375 },
    @{ Sheet = "CRP"; Cell = "I99"; Value = 1157 },
    @{ Sheet = "CRP"; Cell = "J99"; Value = 1398.5 },
    @{ Sheet = "CRP"; Cell = "K99"; Value = 1157 },
    @{ Sheet = "CRP"; Cell = "L99"; Value = 1398.5 },
    @{ Sheet = "CRP"; Cell = "M99"; Value = 341 },
    @{ Sheet = "CRP"; Cell = "N99"; Value = -4394.5 },
    @{ Sheet = "CRP"; Cell = "H105"; Value = 494.6875 },
    @{ Sheet = "CRP"; Cell = "I105"; Value = 518.5833 },
    @{ Sheet = "CRP"; Cell = "J105"; Value = 423 },
    @{ Sheet = "CRP"; Cell = "K105"; Value = 518.5833 },
    @{ Sheet = "CRP"; Cell = "L105"; Value = 423 },
    @{ Sheet = "CRP"; Cell = "M105"; Value = 1228.4167 },
    @{ Sheet = "CRP"; Cell = "N105"; Value = -3917 },
    @{ Sheet = "CRP"; Cell = "H126"; Value = 1217.375 },
    @{ Sheet = "CRP"; Cell = "I126"; Value = 1157 },
    @{ Sheet = "CRP"; Cell = "J126"; Value = 1398.5 },
    @{ Sheet = "CRP"; Cell = "K126"; Value = 3471 },
    @{ Sheet = "CRP"; Cell = "L126"; Value = 4195.5 },
    @{ Sheet = "CRP"; Cell = "M126"; Value = -1001 },
    @{ Sheet = "CRP"; Cell = "N126"; Value = -9135.5 },
    @{ Sheet = "CUL"; Cell = "H8"; Value = 65.818184 },
    @{ Sheet = "CUL"; Cell = "I8"; Value = 65.818184 },
    @{ Sheet = "CUL"; Cell = "K8"; Value = 197.454552 },
    @{ Sheet = "CUL"; Cell = "M8"; Value = -58.45455200000001 },
    @{ Sheet = "CUL"; Cell = "H92"; Value = 350 },
    @{ Sheet = "CUL"; Cell = "I92"; Value = 350 },
    @{ Sheet = "CUL"; Cell = "K92"; Value = 1050 },
    @{ Sheet = "CUL"; Cell = "M92"; Value = 198 },
    @{ Sheet = "CUL"; Cell = "H106"; Value = 3723.1 },
    @{ Sheet = "CUL"; Cell = "J106"; Value = 3723.1 },
    @{ Sheet = "CUL"; Cell = "L106"; Value = 11169.3 },
    @{ Sheet = "CUL"; Cell = "N106"; Value = -13061.3 },
    @{ Sheet = "CUL"; Cell = "H109"; Value = 2247.7778 },
    @{ Sheet = "CUL"; Cell = "I109"; Value = 1680 },
    @{ Sheet = "CUL"; Cell = "K109"; Value = 5040 },
    @{ Sheet = "CUL"; Cell = "M109"; Value = -4000 },
    @{ Sheet = "CUL"; Cell = "H132"; Value = 1274.8235 },
    @{ Sheet = "CUL"; Cell = "I132"; Value = 978.3 },
    @{ Sheet = "CUL"; Cell = "K132"; Value = 8804.699999999999 },
    @{ Sheet = "CUL"; Cell = "M132"; Value = -6274.699999999999 },
    @{ Sheet = "GSM"; Cell = "H70"; Value = 5673.4287 },
    @{ Sheet = "GSM"; Cell = "I70"; Value = 5497.6216 },
    @{ Sheet = "GSM"; Cell = "J70"; Value = 6015.7896 },
    @{ Sheet = "GSM"; Cell = "K70"; Value = 5497.6216 },
    @{ Sheet = "GSM"; Cell = "L70"; Value = 6015.7896 },
    @{ Sheet = "GSM"; Cell = "M70"; Value = -5227.6216 },
    @{ Sheet = "GSM"; Cell = "N70"; Value = -6555.7896 },
    @{ Sheet = "GSM"; Cell = "H73"; Value = 5673.4287 },
    @{ Sheet = "GSM"; Cell = "I73"; Value = 5497.6216 },
    @{ Sheet = "GSM"; Cell = "J73"; Value = 6015.7896 },
    @{ Sheet = "GSM"; Cell = "K73"; Value = 5497.6216 },
    @{ Sheet = "GSM"; Cell = "L73"; Value = 6015.7896 },
    @{ Sheet = "GSM"; Cell = "M73"; Value = -4561.6216 },
    @{ Sheet = "GSM"; Cell = "N73"; Value = -7887.7896 },
    @{ Sheet = "LTW"; Cell = "H7"; Value = 3170.3684 },
    @{ Sheet = "LTW"; Cell = "I7"; Value = 2625.577 },
    @{ Sheet = "LTW"; Cell = "J7"; Value = 4350.75 },
    @{ Sheet = "LTW"; Cell = "K7"; Value = 2625.577 },
    @{ Sheet = "LTW"; Cell = "L7"; Value = 4350.75 },
    @{ Sheet = "LTW"; Cell = "M7"; Value = -2513.577 },
    @{ Sheet = "LTW"; Cell = "N7"; Value = -4574.75 },
    @{ Sheet = "LTW"; Cell = "H22"; Value = 1060.8 },
    @{ Sheet = "LTW"; Cell = "I22"; Value = 500 },
    @{ Sheet = "LTW"; Cell = "J22"; Value = 1201 },
    @{ Sheet = "LTW"; Cell = "K22"; Value = 500 },
    @{ Sheet = "LTW"; Cell = "L22"; Value = 1201 },
    @{ Sheet = "LTW"; Cell = "M22"; Value = -205 },
    @{ Sheet = "LTW"; Cell = "N22"; Value = -1791 },
    @{ Sheet = "LTW"; Cell = "H27"; Value = 1060.8 },
    @{ Sheet = "LTW"; Cell = "I27"; Value = 500 },
    @{ Sheet = "LTW"; Cell = "J27"; Value = 1201 },
    @{ Sheet = "LTW"; Cell = "K27"; Value = 500 },
    @{ Sheet = "LTW"; Cell = "L27"; Value = 1201 },
    @{ Sheet = "LTW"; Cell = "M27"; Value = -393 },
    @{ Sheet = "LTW"; Cell = "N27"; Value = -1415 },
    @{ Sheet = "LTW"; Cell = "H40"; Value = 4006.8125 },
    @{ Sheet = "LTW"; Cell = "I40"; Value = 3662.2307 },
    @{ Sheet = "LTW"; Cell = "K40"; Value = 3662.2307 },
    @{ Sheet = "LTW"; Cell = "M40"; Value = -3526.2307 },
    @{ Sheet = "LTW"; Cell = "H44"; Value = 6000 },
    @{ Sheet = "LTW"; Cell = "J44"; Value = 6000 },
    @{ Sheet = "LTW"; Cell = "L44"; Value = 6000 },
    @{ Sheet = "LTW"; Cell = "M44"; Value = -6912 },
    @{ Sheet = "LTW"; Cell = "H56"; Value = 10896.143 },
    @{ Sheet = "LTW"; Cell = "I56"; Value = 14000 },
    @{ Sheet = "LTW"; Cell = "J56"; Value = 10378.833 },
    @{ Sheet = "LTW"; Cell = "K56"; Value = 14000 },
    @{ Sheet = "LTW"; Cell = "L56"; Value = 10378.833 },
    @{ Sheet = "LTW"; Cell = "M56"; Value = -13309 },
    @{ Sheet = "LTW"; Cell = "N56"; Value = -11760.833 },
    @{ Sheet = "LTW"; Cell = "H61"; Value = 1701501.5 },
    @{ Sheet = "LTW"; Cell = "I61"; Value = 41800.8 },
    @{ Sheet = "LTW"; Cell = "J61"; Value = 10000005 },
    @{ Sheet = "LTW"; Cell = "K61"; Value = 41800.8 },
    @{ Sheet = "LTW"; Cell = "L61"; Value = 10000005 },
    @{ Sheet = "LTW"; Cell = "M61"; Value = -41598.8 },
    @{ Sheet = "LTW"; Cell = "N61"; Value = -10000409 },
    @{ Sheet = "LTW"; Cell = "H113"; Value = 1701501.5 },
    @{ Sheet = "LTW"; Cell = "I113"; Value = 41800.8 },
    @{ Sheet = "LTW"; Cell = "J113"; Value = 10000005 },
    @{ Sheet = "LTW"; Cell = "K113"; Value = 41800.8 },
    @{ Sheet = "LTW"; Cell = "L113"; Value = 10000005 },
    @{ Sheet = "LTW"; Cell = "M113"; Value = -39630.8 },
    @{ Sheet = "LTW"; Cell = "N113"; Value = -10004345 },
    @{ Sheet = "LTW"; Cell = "H126"; Value = 3170.3684 },
    @{ Sheet = "LTW"; Cell = "I126"; Value = 2625.577 },
    @{ Sheet = "LTW"; Cell = "J126"; Value = 4350.75 },
    @{ Sheet = "LTW"; Cell = "K126"; Value = 7876.731000000001 },
    @{ Sheet = "LTW"; Cell = "L126"; Value = 13052.25 },
    @{ Sheet = "LTW"; Cell = "M126"; Value = -5406.731000000001 },
    @{ Sheet = "LTW"; Cell = "N126"; Value = -17992.25 },
    @{ Sheet = "WVR"; Cell = "H113"; Value = 587.6896400000001 },
    @{ Sheet = "WVR"; Cell = "I113"; Value = 392.35294 },
    @{ Sheet = "WVR"; Cell = "J113"; Value = 864.4167 },
    @{ Sheet = "WVR"; Cell = "K113"; Value = 1177.05882 },
    @{ Sheet = "WVR"; Cell = "L113"; Value = 2593.2501 },
    @{ Sheet = "WVR"; Cell = "M113"; Value = 992.94118 },
    @{ Sheet = "WVR"; Cell = "N113"; Value = -6933.2501 },
    @{ Sheet = "WVR"; Cell = "H122"; Value = 1000 },
    @{ Sheet = "WVR"; Cell = "I122"; Value = 1000 },
    @{ Sheet = "WVR"; Cell = "K122"; Value = 3000 },
    @{ Sheet = "WVR"; Cell = "M122"; Value = -550 }
)

foreach ($e in $edits) {
    $ws = $wb.Worksheets.Item($e.Sheet)
    $ws.Range($e.Cell).Value = $e.Value
}

Write-Host "Applied" $edits.Length "cell updates."
